$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (column labels; text unchanged from before, rewritten for completeness)
$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "Prognose"
$ws.Range("C1").Value = "surveys"
$ws.Range("D1").Value = "production"
$ws.Range("E1").Value = "orders"
$ws.Range("F1").Value = "turnover"
$ws.Range("G1").Value = "financial"
$ws.Range("H1").Value = "labor market"
$ws.Range("I1").Value = "prices"
$ws.Range("J1").Value = "national accounts"
$ws.Range("K1").Value = "Revision"

# Row 2: 2025-03-30
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-03-30"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = [double]"0.30995030401720053"
$ws.Range("C2").Value = [double]"0"
$ws.Range("D2").Value = [double]"0"
$ws.Range("E2").Value = [double]"0"
$ws.Range("F2").Value = [double]"0"
$ws.Range("G2").Value = [double]"0"
$ws.Range("H2").Value = [double]"0"
$ws.Range("I2").Value = [double]"0"
$ws.Range("J2").Value = [double]"0"
$ws.Range("K2").Value = [double]"0"

# Row 3: 2025-04-15
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-04-15"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = [double]"0.3040718621835608"
$ws.Range("C3").Value = [double]"0"
$ws.Range("D3").Value = [double]"-0.0053712655749828921"
$ws.Range("E3").Value = [double]"0.00021751932506205947"
$ws.Range("F3").Value = [double]"-0.002154772126341054"
$ws.Range("G3").Value = [double]"0.00034124346655849336"
$ws.Range("H3").Value = [double]"-0.00025933402887663059"
$ws.Range("I3").Value = [double]"-9.5776351640067405e-05"
$ws.Range("J3").Value = [double]"0"
$ws.Range("K3").Value = [double]"0.0014439434565803877"

# Row 4: 2025-04-30
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2025-04-30"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = [double]"0.30095461904884196"
$ws.Range("C4").Value = [double]"-0.0027007833276229574"
$ws.Range("D4").Value = [double]"0"
$ws.Range("E4").Value = [double]"0.0010560876359296304"
$ws.Range("F4").Value = [double]"6.4253866815975968e-05"
$ws.Range("G4").Value = [double]"0"
$ws.Range("H4").Value = [double]"0.00015706546942532201"
$ws.Range("I4").Value = [double]"-0.0018380868719951263"
$ws.Range("J4").Value = [double]"0.00018165535839600343"
$ws.Range("K4").Value = [double]"-3.7435265667684625e-05"

# Row 5: 2025-05-15
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2025-05-15"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = [double]"0.30585880894429501"
$ws.Range("C5").Value = [double]"0.011442692473081601"
$ws.Range("D5").Value = [double]"-0.0057495412481816622"
$ws.Range("E5").Value = [double]"0.0007151817618350388"
$ws.Range("F5").Value = [double]"0.0011197859938279486"
$ws.Range("G5").Value = [double]"-0.0020268756620531016"
$ws.Range("H5").Value = [double]"9.7950226966335153e-05"
$ws.Range("I5").Value = [double]"-0.00055164293768455877"
$ws.Range("J5").Value = [double]"0"
$ws.Range("K5").Value = [double]"-0.00014336071233855829"

# Row 6: 2025-05-30
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2025-05-30"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = [double]"0.32274297855076395"
$ws.Range("C6").Value = [double]"0.019974638536148268"
$ws.Range("D6").Value = [double]"0"
$ws.Range("E6").Value = [double]"-0.00068213905944316232"
$ws.Range("F6").Value = [double]"3.7157577308440897e-05"
$ws.Range("G6").Value = [double]"0"
$ws.Range("H6").Value = [double]"-2.4256554430993207e-06"
$ws.Range("I6").Value = [double]"-0.0012136065679540632"
$ws.Range("J6").Value = [double]"0"
$ws.Range("K6").Value = [double]"-0.0012294552241474133"

# Row 7: 2025-06-15
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2025-06-15"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = [double]"0.30761177401231937"
$ws.Range("C7").Value = [double]"0"
$ws.Range("D7").Value = [double]"-0.0044273649000430996"
$ws.Range("E7").Value = [double]"-0.002600686046821983"
$ws.Range("F7").Value = [double]"-0.0085691464355290448"
$ws.Range("G7").Value = [double]"0.0011223611278006606"
$ws.Range("H7").Value = [double]"0"
$ws.Range("I7").Value = [double]"0.00024263536285114073"
$ws.Range("J7").Value = [double]"0"
$ws.Range("K7").Value = [double]"-0.00089900364670225663"

# Row 8: 2025-06-30
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "2025-06-30"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = [double]"0.25019645202814311"
$ws.Range("C8").Value = [double]"-0.05722038976859592"
$ws.Range("D8").Value = [double]"0"
$ws.Range("E8").Value = [double]"-6.8611443609356074e-05"
$ws.Range("F8").Value = [double]"-0.00027537766601826353"
$ws.Range("G8").Value = [double]"0"
$ws.Range("H8").Value = [double]"-6.1250545798344249e-05"
$ws.Range("I8").Value = [double]"-4.657826615698986e-05"
$ws.Range("J8").Value = [double]"0"
$ws.Range("K8").Value = [double]"0.00025688570600274074"

# Row 9: 2025-07-15
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "2025-07-15"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = [double]"0.25809068282547526"
$ws.Range("C9").Value = [double]"0"
$ws.Range("D9").Value = [double]"0.01008181288044566"
$ws.Range("E9").Value = [double]"-0.0026426501251446276"
$ws.Range("F9").Value = [double]"-0.0033187827016197525"
$ws.Range("G9").Value = [double]"0.0030302215888688552"
$ws.Range("H9").Value = [double]"-0.00014405567747046057"
$ws.Range("I9").Value = [double]"0.00042984276341204684"
$ws.Range("J9").Value = [double]"0"
$ws.Range("K9").Value = [double]"0.00045784206884036394"

# Row 10: 2025-07-30
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "2025-07-30"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = [double]"0.35007908033238705"
$ws.Range("C10").Value = [double]"0.093529143852494098"
$ws.Range("D10").Value = [double]"0"
$ws.Range("E10").Value = [double]"-0.001166278204748407"
$ws.Range("F10").Value = [double]"-2.3447854393231712e-05"
$ws.Range("G10").Value = [double]"0"
$ws.Range("H10").Value = [double]"2.6183822140070506e-05"
$ws.Range("I10").Value = [double]"0.0011399283016942101"
$ws.Range("J10").Value = [double]"-0.0026587251606346313"
$ws.Range("K10").Value = [double]"0.0011415927503597989"

# Row 11: 2025-08-15
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "2025-08-15"
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").Value = [double]"0.33758961689237632"
$ws.Range("C11").Value = [double]"0"
$ws.Range("D11").Value = [double]"-0.043967231907637586"
$ws.Range("E11").Value = [double]"0.0058923252672639604"
$ws.Range("F11").Value = [double]"0.0076012814540390976"
$ws.Range("G11").Value = [double]"0.0058231723950838638"
$ws.Range("H11").Value = [double]"0.0016880113072155977"
$ws.Range("I11").Value = [double]"0.0074116738373836127"
$ws.Range("J11").Value = [double]"0"
$ws.Range("K11").Value = [double]"0.003061304206640747"

# Row 12: 2025-08-30
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "2025-08-30"
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").Value = [double]"0.26823028143095728"
$ws.Range("C12").Value = [double]"-0.064158193642574288"
$ws.Range("D12").Value = [double]"0"
$ws.Range("E12").Value = [double]"0.0028090780366104404"
$ws.Range("F12").Value = [double]"2.6474738809491305e-05"
$ws.Range("G12").Value = [double]"0"
$ws.Range("H12").Value = [double]"-1.4603632978030841e-05"
$ws.Range("I12").Value = [double]"-0.0038343266517918131"
$ws.Range("J12").Value = [double]"0"
$ws.Range("K12").Value = [double]"-0.0041877643094948569"
